$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COOLING")

$ws.Columns("C:E").Insert()
$ws.Columns("C:E").ColumnWidth = 20.16666666666667

$ws.Range("C1").Value = "primary_components"
$ws.Range("D1").Value = "secondary_components"
$ws.Range("E1").Value = "tertiary_components"

$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"

$ws.Range("C3").Value = "CH2"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "CT2"

$ws.Range("C4").Value = "CH2"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "CT1"

$ws.Range("C5").Value = "CH1"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "HEX1"

$ws.Range("C6").Value = "CH1"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "CT1"

$ws.Range("C7").Value = "AC1"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"

$wsHeating = $wb.Worksheets.Item("HEATING")
$wsHeating.Range("C2:E2").Copy()
$ws.Range("C2:E7").PasteSpecial(-4122)

$wsHeating.Range("C1:E1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

$ws.Activate()
$ws.Range("A14").Select()
$ws.Range("D18").Select()
